$wb = $excel.ActiveWorkbook

# ---- Sheet: Главные ----
$ws = $wb.Worksheets.Item("Главные")

$ws.Range("C3").Value = 34
$ws.Range("D3").Value = 571
$ws.Range("E3").Value = 259
$ws.Range("F3").Value = 312
$ws.Range("G3").Value = 16.79
$ws.Range("H3").Value = 7.62
$ws.Range("I3").Value = 9.18
$ws.Range("J3").Value = 127
$ws.Range("K3").Value = 131
$ws.Range("W3").Value = 14
$ws.Range("C5").Value = 34
$ws.Range("D5").Value = 531
$ws.Range("E5").Value = 266
$ws.Range("F5").Value = 265
$ws.Range("G5").Value = 15.62
$ws.Range("H5").Value = 7.82
$ws.Range("I5").Value = 7.79
$ws.Range("J5").Value = 128
$ws.Range("K5").Value = 125
$ws.Range("C18").Value = 33
$ws.Range("D18").Value = 553
$ws.Range("E18").Value = 270
$ws.Range("G18").Value = 16.76
$ws.Range("H18").Value = 8.18
$ws.Range("I18").Value = 8.58
$ws.Range("J18").Value = 105
$ws.Range("C19").Value = 26
$ws.Range("D19").Value = 438
$ws.Range("E19").Value = 220
$ws.Range("G19").Value = 16.85
$ws.Range("H19").Value = 8.460000000000001
$ws.Range("I19").Value = 8.380000000000001
$ws.Range("J19").Value = 105
$ws.Range("C21").Value = 32
$ws.Range("D21").Value = 468
$ws.Range("E21").Value = 214
$ws.Range("F21").Value = 254
$ws.Range("G21").Value = 14.63
$ws.Range("H21").Value = 6.69
$ws.Range("I21").Value = 7.94
$ws.Range("J21").Value = 97
$ws.Range("K21").Value = 112
$ws.Range("C23").Value = 21
$ws.Range("D23").Value = 282
$ws.Range("E23").Value = 116
$ws.Range("F23").Value = 166
$ws.Range("G23").Value = 13.43
$ws.Range("H23").Value = 5.52
$ws.Range("I23").Value = 7.9
$ws.Range("J23").Value = 53
$ws.Range("K23").Value = 68
$ws.Range("W23").Value = 6

# Update as_of_utc timestamps for rows 2-26
for ($r = 2; $r -le 26; $r++) {
    $ws.Cells.Item($r, 27).Value = "2025-12-20 03:08:47"
}

# ---- Sheet: Линейные ----
$ws = $wb.Worksheets.Item("Линейные")

$ws.Range("C4").Value = 18
$ws.Range("D4").Value = 254
$ws.Range("E4").Value = 122
$ws.Range("G4").Value = 14.11
$ws.Range("H4").Value = 6.78
$ws.Range("I4").Value = 7.33
$ws.Range("J4").Value = 61
$ws.Range("C7").Value = 21
$ws.Range("D7").Value = 387
$ws.Range("E7").Value = 166
$ws.Range("G7").Value = 18.43
$ws.Range("I7").Value = 10.52
$ws.Range("J7").Value = 63
$ws.Range("C8").Value = 30
$ws.Range("D8").Value = 481
$ws.Range("E8").Value = 203
$ws.Range("F8").Value = 278
$ws.Range("G8").Value = 16.03
$ws.Range("H8").Value = 6.77
$ws.Range("I8").Value = 9.27
$ws.Range("J8").Value = 84
$ws.Range("K8").Value = 109
$ws.Range("W8").Value = 6
$ws.Range("C12").Value = 30
$ws.Range("D12").Value = 500
$ws.Range("E12").Value = 240
$ws.Range("F12").Value = 260
$ws.Range("G12").Value = 16.67
$ws.Range("H12").Value = 8
$ws.Range("I12").Value = 8.67
$ws.Range("J12").Value = 110
$ws.Range("K12").Value = 120
$ws.Range("W12").Value = 12
$ws.Range("C18").Value = 35
$ws.Range("D18").Value = 594
$ws.Range("E18").Value = 279
$ws.Range("F18").Value = 315
$ws.Range("G18").Value = 16.97
$ws.Range("H18").Value = 7.97
$ws.Range("I18").Value = 9
$ws.Range("J18").Value = 132
$ws.Range("K18").Value = 135
$ws.Range("C26").Value = 31
$ws.Range("D26").Value = 606
$ws.Range("E26").Value = 274
$ws.Range("F26").Value = 332
$ws.Range("G26").Value = 19.55
$ws.Range("H26").Value = 8.84
$ws.Range("I26").Value = 10.71
$ws.Range("J26").Value = 112
$ws.Range("K26").Value = 106

# Update as_of_utc timestamps for rows 2-26
for ($r = 2; $r -le 26; $r++) {
    $ws.Cells.Item($r, 27).Value = "2025-12-20 03:08:47"
}

